# Weekly update: insert a new week's worth of "Frutilla" price data
# (date serial 45041) at the top of the historical block, pushing the
# existing rows down by 3 (one row per quality grade).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows right before row 1068 (shifts old 1068:1125 -> 1071:1128)
$ws.Range("A1068:T1070").Insert()

# Columns that are constant across every data row in this sheet/series.
$mercadoId = 8
$mercado   = "Terminal La Palmera de La Serena"
$region    = "Coquimbo"
$codreg    = 4
$tipo      = "Fruta"
$productoId = 100101
$producto   = "Berries"
$categoriaId = 100112025
$categoria   = "Frutilla"
$variedad    = "Sin especificar"
$unidad      = "$/bandeja 7 kilos"
$origen      = "Provincia de Melipilla"
$kgUnidad    = 7
$fecha       = 45041

# New row 1068: Calidad = Especial
$ws.Cells.Item(1068, 1).Value = $mercadoId
$ws.Cells.Item(1068, 2).Value = $mercado
$ws.Cells.Item(1068, 3).Value = $region
$ws.Cells.Item(1068, 4).Value = $fecha
$ws.Cells.Item(1068, 5).Value = $codreg
$ws.Cells.Item(1068, 6).Value = $tipo
$ws.Cells.Item(1068, 7).Value = $productoId
$ws.Cells.Item(1068, 8).Value = $producto
$ws.Cells.Item(1068, 9).Value = $categoriaId
$ws.Cells.Item(1068, 10).Value = $categoria
$ws.Cells.Item(1068, 11).Value = $variedad
$ws.Cells.Item(1068, 12).Value = "Especial"
$ws.Cells.Item(1068, 13).Value = 400
$ws.Cells.Item(1068, 14).Value = 15000
$ws.Cells.Item(1068, 15).Value = 16000
$ws.Cells.Item(1068, 16).Value = 15500
$ws.Cells.Item(1068, 17).Value = $unidad
$ws.Cells.Item(1068, 18).Value = $origen
$ws.Cells.Item(1068, 19).Value = 2214
$ws.Cells.Item(1068, 20).Value = $kgUnidad

# New row 1069: Calidad = Primera
$ws.Cells.Item(1069, 1).Value = $mercadoId
$ws.Cells.Item(1069, 2).Value = $mercado
$ws.Cells.Item(1069, 3).Value = $region
$ws.Cells.Item(1069, 4).Value = $fecha
$ws.Cells.Item(1069, 5).Value = $codreg
$ws.Cells.Item(1069, 6).Value = $tipo
$ws.Cells.Item(1069, 7).Value = $productoId
$ws.Cells.Item(1069, 8).Value = $producto
$ws.Cells.Item(1069, 9).Value = $categoriaId
$ws.Cells.Item(1069, 10).Value = $categoria
$ws.Cells.Item(1069, 11).Value = $variedad
$ws.Cells.Item(1069, 12).Value = "Primera"
$ws.Cells.Item(1069, 13).Value = 300
$ws.Cells.Item(1069, 14).Value = 13000
$ws.Cells.Item(1069, 15).Value = 14000
$ws.Cells.Item(1069, 16).Value = 13500
$ws.Cells.Item(1069, 17).Value = $unidad
$ws.Cells.Item(1069, 18).Value = $origen
$ws.Cells.Item(1069, 19).Value = 1929
$ws.Cells.Item(1069, 20).Value = $kgUnidad

# New row 1070: Calidad = Segunda
$ws.Cells.Item(1070, 1).Value = $mercadoId
$ws.Cells.Item(1070, 2).Value = $mercado
$ws.Cells.Item(1070, 3).Value = $region
$ws.Cells.Item(1070, 4).Value = $fecha
$ws.Cells.Item(1070, 5).Value = $codreg
$ws.Cells.Item(1070, 6).Value = $tipo
$ws.Cells.Item(1070, 7).Value = $productoId
$ws.Cells.Item(1070, 8).Value = $producto
$ws.Cells.Item(1070, 9).Value = $categoriaId
$ws.Cells.Item(1070, 10).Value = $categoria
$ws.Cells.Item(1070, 11).Value = $variedad
$ws.Cells.Item(1070, 12).Value = "Segunda"
$ws.Cells.Item(1070, 13).Value = 240
$ws.Cells.Item(1070, 14).Value = 10000
$ws.Cells.Item(1070, 15).Value = 11000
$ws.Cells.Item(1070, 16).Value = 10500
$ws.Cells.Item(1070, 17).Value = $unidad
$ws.Cells.Item(1070, 18).Value = $origen
$ws.Cells.Item(1070, 19).Value = 1500
$ws.Cells.Item(1070, 20).Value = $kgUnidad
